$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - same style as other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Timestamps for F2:F21 - store as plain text to match inlineStr type in source
$timestamps = @(
    "2021-10-05 10:50:38.885373",
    "2021-10-05 10:50:38.885387",
    "2021-10-05 10:50:38.885391",
    "2021-10-05 10:50:38.885395",
    "2021-10-05 10:50:38.885398",
    "2021-10-05 10:50:38.885401",
    "2021-10-05 10:50:38.885404",
    "2021-10-05 10:50:38.885407",
    "2021-10-05 10:50:38.885411",
    "2021-10-05 10:50:38.885414",
    "2021-10-05 10:50:38.885417",
    "2021-10-05 10:50:38.885420",
    "2021-10-05 10:50:38.885423",
    "2021-10-05 10:50:38.885427",
    "2021-10-05 10:50:38.885430",
    "2021-10-05 10:50:38.885433",
    "2021-10-05 10:50:38.885437",
    "2021-10-05 10:50:38.885440",
    "2021-10-05 10:50:38.885443",
    "2021-10-05 10:50:38.885446"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
